# Natmi following Dr Hou advice
# Update existing rows 2-5 with recomputed statistics, and insert two new
# target-cluster rows (M1, M2) between the existing "ECs"/"FAPs" rows and
# the "Neutro"/"sCs" rows, shifting the latter down to rows 6-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (Target cluster: ECs) - updated values
# ---------------------------------------------------------------------
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.619088000000001
$ws.Range("H2").Value = 13.857264
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 2.569008
$ws.Range("N2").Value = 5.138016
$ws.Range("O2").Value = 0.2577238367773512
$ws.Range("P2").Value = 0.2133397199412101
$ws.Range("Q2").Value = 11.866474024704
$ws.Range("R2").Value = 71.19884414822401
$ws.Range("S2").Value = 0.2577238367773512
$ws.Range("T2").Value = 0.2133397199412101

# ---------------------------------------------------------------------
# Row 3 (Target cluster: FAPs) - updated values
# ---------------------------------------------------------------------
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.619088000000001
$ws.Range("H3").Value = 13.857264
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.043481333333333
$ws.Range("N3").Value = 6.130444
$ws.Range("O3").Value = 0.2050028063787906
$ws.Range("P3").Value = 0.2545471259870097
$ws.Range("Q3").Value = 9.439020105024001
$ws.Range("R3").Value = 84.951180945216
$ws.Range("S3").Value = 0.2050028063787906
$ws.Range("T3").Value = 0.2545471259870097

# ---------------------------------------------------------------------
# Row 4 was "Neutro" -> becomes "M1" with recomputed values, and moves to
# row 6 later; row 4 now holds the new "M1" data.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.619088000000001
$ws.Range("H4").Value = 13.857264
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01206
$ws.Range("N4").Value = 0.03618
$ws.Range("O4").Value = 0.001209863679496076
$ws.Range("P4").Value = 0.001502259056311421
$ws.Range("Q4").Value = 0.05570620128000001
$ws.Range("R4").Value = 0.50135581152
$ws.Range("S4").Value = 0.001209863679496076
$ws.Range("T4").Value = 0.001502259056311421

# ---------------------------------------------------------------------
# Row 5 was "sCs" -> becomes "M2" with recomputed values.
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd8"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.619088000000001
$ws.Range("H5").Value = 13.857264
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03685233333333333
$ws.Range("N5").Value = 0.110557
$ws.Range("O5").Value = 0.003697039768215801
$ws.Range("P5").Value = 0.004590526658060303
$ws.Range("Q5").Value = 0.170224170672
$ws.Range("R5").Value = 1.532017536048
$ws.Range("S5").Value = 0.003697039768215801
$ws.Range("T5").Value = 0.004590526658060303

# ---------------------------------------------------------------------
# Row 6 (new): Target cluster Neutro - the original row-4 data, recomputed.
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd8"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.619088000000001
$ws.Range("H6").Value = 13.857264
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.055205
$ws.Range("N6").Value = 6.165615
$ws.Range("O6").Value = 0.206178928973361
$ws.Range("P6").Value = 0.2560074895378535
$ws.Range("Q6").Value = 9.493172753040001
$ws.Range("R6").Value = 85.43855477736
$ws.Range("S6").Value = 0.206178928973361
$ws.Range("T6").Value = 0.2560074895378535

# ---------------------------------------------------------------------
# Row 7 (new): Target cluster sCs - the original row-5 data, recomputed.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd8"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.619088000000001
$ws.Range("H7").Value = 13.857264
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.2514585
$ws.Range("N7").Value = 6.502917
$ws.Range("O7").Value = 0.3261875244227855
$ws.Range("P7").Value = 0.2700128788195549
$ws.Range("Q7").Value = 15.018772939848
$ws.Range("R7").Value = 90.11263763908801
$ws.Range("S7").Value = 0.3261875244227855
$ws.Range("T7").Value = 0.2700128788195549
